# Refresh market-price derived columns (currentAveragePrice.. LeveProfitHQ, H:N)
# on the Gungnir_Profits leve-crafting sheets. Values below come from the
# scheduled market-data sync; only H:N are touched, A:G (leve metadata) are
# left untouched.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 107: Another Man's Ink
$ws.Range("H107").Value = 266
$ws.Range("I107").Value = 230
$ws.Range("J107").Value = 350
$ws.Range("K107").Value = 230
$ws.Range("L107").Value = 350
$ws.Range("M107").Value = 1690
$ws.Range("N107").Value = -4190

# Row 137: Cutting Edge of Culinary Quality
$ws.Range("H137").Value = 2299.9092
$ws.Range("I137").Value = 1642.7142
$ws.Range("J137").Value = 3450
$ws.Range("K137").Value = 4928.142599999999
$ws.Range("L137").Value = 10350
$ws.Range("M137").Value = -2378.142599999999
$ws.Range("N137").Value = -15450

$ws = $wb.Worksheets.Item("ARM")
# Row 2: Ain't Got No Ingots
$ws.Range("H2").Value = 849.8125
$ws.Range("I2").Value = 799.75
$ws.Range("K2").Value = 799.75
$ws.Range("M2").Value = -686.75

# Row 7: Distill It Yourself
$ws.Range("H7").Value = 41000
$ws.Range("J7").Value = 41000
$ws.Range("L7").Value = 41000
$ws.Range("N7").Value = -41228

# Row 32: Ingot We Trust
$ws.Range("H32").Value = 32286426
$ws.Range("I32").Value = 55570300
$ws.Range("J32").Value = 47216.46
$ws.Range("K32").Value = 55570300
$ws.Range("L32").Value = 47216.46
$ws.Range("M32").Value = -55570013
$ws.Range("N32").Value = -47790.46

# Row 110: Scheduled Maintenance
$ws.Range("H110").Value = 1446.75
$ws.Range("I110").Value = 888.0833
$ws.Range("J110").Value = 3122.75
$ws.Range("K110").Value = 888.0833
$ws.Range("L110").Value = 3122.75
$ws.Range("M110").Value = 1156.9167
$ws.Range("N110").Value = -7212.75

# Row 116: No Scope
$ws.Range("H116").Value = 849.8125
$ws.Range("I116").Value = 799.75
$ws.Range("K116").Value = 799.75
$ws.Range("M116").Value = 1494.25

# Row 117: Signed, Shield, Delivered
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()

# Row 122: Haste for High Durium
$ws.Range("H122").Value = 9380
$ws.Range("I122").Value = 11506.667
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 34520.001
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -32070.001
$ws.Range("N122").Value = -13900

$ws = $wb.Worksheets.Item("BSM")
# Row 3: Hells Bells
$ws.Range("H3").Value = 849.8125
$ws.Range("I3").Value = 799.75
$ws.Range("K3").Value = 799.75
$ws.Range("M3").Value = -685.75

# Row 7: Thank You for Your Business
$ws.Range("H7").Value = 6752625
$ws.Range("I7").Value = 9000833
$ws.Range("J7").Value = 8000
$ws.Range("K7").Value = 9000833
$ws.Range("L7").Value = 8000
$ws.Range("M7").Value = -9000720
$ws.Range("N7").Value = -8226

$ws = $wb.Worksheets.Item("CRP")
# Row 23: Nothing to Hide
$ws.Range("H23").Value = 13500
$ws.Range("I23").Value = 13500
$ws.Range("K23").Value = 13500
$ws.Range("M23").Value = -13260

# Row 27: Behind the Mask
$ws.Range("H27").Value = 13500
$ws.Range("I27").Value = 13500
$ws.Range("K27").Value = 13500
$ws.Range("M27").Value = -13308

# Row 31: Wall Not Found
$ws.Range("H31").Value = 3418.1924
$ws.Range("I31").Value = 2639.5881
$ws.Range("J31").Value = 4888.8887
$ws.Range("K31").Value = 2639.5881
$ws.Range("L31").Value = 4888.8887
$ws.Range("M31").Value = -2344.5881
$ws.Range("N31").Value = -5478.8887

# Row 32: Daddy's Little Girl
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("M32").ClearContents()
$ws.Range("N32").ClearContents()

# Row 34: Armoires of the Rich and Famous
$ws.Range("H34").Value = 3418.1924
$ws.Range("I34").Value = 2639.5881
$ws.Range("J34").Value = 4888.8887
$ws.Range("K34").Value = 2639.5881
$ws.Range("L34").Value = 4888.8887
$ws.Range("M34").Value = -2437.5881
$ws.Range("N34").Value = -5292.8887

# Row 62: Splinter in the Sewers
$ws.Range("H62").Value = 3582.3044
$ws.Range("I62").Value = 2721.2727
$ws.Range("J62").Value = 4371.5835
$ws.Range("K62").Value = 2721.2727
$ws.Range("L62").Value = 4371.5835
$ws.Range("M62").Value = -2097.2727
$ws.Range("N62").Value = -5619.5835

# Row 65: The Lumber of Their Discontent (L)
$ws.Range("H65").Value = 3582.3044
$ws.Range("I65").Value = 2721.2727
$ws.Range("J65").Value = 4371.5835
$ws.Range("K65").Value = 13606.3635
$ws.Range("L65").Value = 21857.9175
$ws.Range("M65").Value = -10486.3635
$ws.Range("N65").Value = -28097.9175

$ws = $wb.Worksheets.Item("CUL")
# Row 70: Persona non Gratin
$ws.Range("H70").Value = 700
$ws.Range("I70").Value = 700
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 2100
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -1785
$ws.Range("N70").ClearContents()

# Row 73: Recipe for Disaster (L)
$ws.Range("H73").Value = 700
$ws.Range("I73").Value = 700
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 2100
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -1008
$ws.Range("N73").ClearContents()

# Row 122: Salt of the North
$ws.Range("H122").Value = 33340804
$ws.Range("I122").Value = 62500396
$ws.Range("J122").Value = 15557
$ws.Range("K122").Value = 562503564
$ws.Range("L122").Value = 140013
$ws.Range("M122").Value = -562501114
$ws.Range("N122").Value = -144913

# Row 131: The Mountain Steeped
$ws.Range("H131").Value = 862.53
$ws.Range("J131").Value = 867.4796
$ws.Range("L131").Value = 2602.4388
$ws.Range("N131").Value = -12682.4388

$ws = $wb.Worksheets.Item("GSM")
# Row 132: On Board for Lar
$ws.Range("H132").Value = 16399.6
$ws.Range("I132").Value = 3990
$ws.Range("J132").Value = 19502
$ws.Range("K132").Value = 11970
$ws.Range("L132").Value = 58506
$ws.Range("M132").Value = -9440
$ws.Range("N132").Value = -63566

$ws = $wb.Worksheets.Item("LTW")
# Row 46: Supply Side Logic
$ws.Range("H46").Value = 3788342.8
$ws.Range("I46").Value = 5208788
$ws.Range("J46").Value = 488.33334
$ws.Range("K46").Value = 5208788
$ws.Range("L46").Value = 488.33334
$ws.Range("M46").Value = -5208600
$ws.Range("N46").Value = -864.33334

# Row 111: Glove Me Tender
$ws.Range("H111").Value = 40300
$ws.Range("J111").Value = 40300
$ws.Range("L111").Value = 40300

# Row 136: Respect for Br'aax
$ws.Range("H136").Value = 3503.25
$ws.Range("I136").Value = 2713.1667
$ws.Range("J136").Value = 5873.5
$ws.Range("K136").Value = 8139.500100000001
$ws.Range("L136").Value = 17620.5
$ws.Range("M136").Value = -5589.500100000001
$ws.Range("N136").Value = -22720.5

$ws = $wb.Worksheets.Item("WVR")
# Row 81: Where the Dragonflies, the Net Catches
$ws.Range("H81").Value = 50000830
$ws.Range("I81").Value = 55556370
$ws.Range("J81").Value = 1000
$ws.Range("K81").Value = 111112740
$ws.Range("L81").Value = 2000
$ws.Range("M81").Value = -111111679
$ws.Range("N81").Value = -4122

# Row 84: To Kill a Dragon on Nameday (L)
$ws.Range("H84").Value = 55556370
$ws.Range("I84").Value = 55556370
$ws.Range("J84").Value = 1000
$ws.Range("K84").Value = 555563700
$ws.Range("L84").Value = 10000
$ws.Range("M84").Value = -10000
